# XML Update 5 - FIB, MCQ, FIB/MCQ Complete
#
# 1. Update the stored selection on "Sheet4" (E18 -> A7:A10) without
#    disturbing which sheet/tab ends up active.
# 2. Duplicate "Sheet6" twice, producing "Sheet6 (2)" and "Sheet6 (3)".
# 3. Fill in "Sheet6 (3)" (FIB/MCQ + FIBTest) first, then "Sheet6 (2)"
#    (FIB/MCQ + aaaMCQTest) -- this ordering reproduces the shared-string
#    insertion order seen in the target workbook.
# 4. Leave "Sheet6 (2)" as the active sheet/tab, matching the saved file.

$wb = $excel.ActiveWorkbook

# --- Step 1: update the remembered selection on Sheet4 -------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("A7:A10").Select()

# --- Step 2: duplicate Sheet6 twice ---------------------------------------
$srcSheet = $wb.Worksheets.Item("Sheet6")

$srcSheet.Copy($null, $srcSheet)
$ws2 = $wb.Worksheets.Item("Sheet6 (2)")

$srcSheet.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item("Sheet6 (3)")

# --- Step 3a: fill in Sheet6 (3) first ------------------------------------
$ws3.Range("C2").Value = "FIB/MCQ"
$ws3.Range("B4").Value = "FIBTest"
$ws3.Range("A7").Value = "C1"
$ws3.Range("A8").Value = "C2"
$ws3.Range("A9").Value = "C3"
$ws3.Range("A10").Value = "C4"
$ws3.Range("A7:A10").Style = "Normal"
$ws3.Range("D12").Select()

# --- Step 3b: fill in Sheet6 (2) second -----------------------------------
$ws2.Range("C2").Value = "FIB/MCQ"
$ws2.Range("B4").Value = "aaaMCQTest"
$ws2.Range("A7").Value = "C1"
$ws2.Range("A8").Value = "C2"
$ws2.Range("A9").Value = "C3"
$ws2.Range("A10").Value = "C4"
$ws2.Range("A7:A10").Style = "Normal"
$ws2.Range("B4").Select()

# --- Step 4: make Sheet6 (2) the active tab -------------------------------
$ws2.Activate()
